$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping it stored as text,
# and without leaving behind a lingering Text number-format/style on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Column D price updates (stored as text strings in the sheet)
Set-TextValue "D2"  "283.43"
Set-TextValue "D3"  "20.61"
Set-TextValue "D4"  "6.220"
Set-TextValue "D5"  "0.06170"
Set-TextValue "D6"  "3.585"
Set-TextValue "D7"  "6.545"
Set-TextValue "D8"  "1.503"
Set-TextValue "D9"  "0.8197"
Set-TextValue "D10" "0.01381"
Set-TextValue "D11" "0.1628"
Set-TextValue "D12" "0.08424"
Set-TextValue "D13" "0.03467"
Set-TextValue "D14" "0.03214"
Set-TextValue "D15" "0.09184"
Set-TextValue "D16" "3.715"
Set-TextValue "D17" "0.001642"
Set-TextValue "D18" "0.04711"
Set-TextValue "D19" "0.006457"
Set-TextValue "D20" "0.006173"
Set-TextValue "D22" "0.0001603"
Set-TextValue "D23" "3.844"
Set-TextValue "D40" "0.04717"
Set-TextValue "D41" "0.007224"

# Rows 42 and 43: BKEXToken / CEJI swapped positions, with new price data
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1099"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003567"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.01147"
Set-TextValue "D45" "0.00006863"
Set-TextValue "D47" "1.102"
Set-TextValue "D48" "0.002849"
Set-TextValue "D49" "0.00001904"
Set-TextValue "D50" "0.01242"
